# Refresh crypto price (D) and 1h volume change (E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'51.500.02"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.09%  "
$c = $ws.Range("D3")
$c.Value = "'3.055.50"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +2.08%  "
$c = $ws.Range("D4")
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "
$c = $ws.Range("D5")
$c.Value = "'385.48"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.93%  "
$c = $ws.Range("D6")
$c.Value = "'103.24"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("E7").Value = "  -0.61%  "
$c = $ws.Range("D9")
$c.Value = "'0.586"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.70%  "
$c = $ws.Range("D10")
$c.Value = "'36.81"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.15%  "
$ws.Range("E11").Value = "  -0.01%  "
$ws.Range("E12").Value = "  -0.13%  "
$c = $ws.Range("D13")
$c.Value = "'3.536.78"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +1.96%  "
$c = $ws.Range("D14")
$c.Value = "'18.56"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +0.53%  "
$c = $ws.Range("D15")
$c.Value = "'7.77"
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.Value = "'3.054.41"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +1.86%  "
$c = $ws.Range("D17")
$c.Value = "'0.974"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -2.86%  "
$c = $ws.Range("D18")
$c.Value = "'10.69"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -5.39%  "
$c = $ws.Range("D19")
$c.Value = "'51.581.69"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.04%  "
$c = $ws.Range("D20")
$c.Value = "'3.16"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +0.28%  "
$c = $ws.Range("D21")
$c.Value = "'12.45"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.37%  "
$c = $ws.Range("D22")
$c.Value = "'0.0₃0964"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.02%  "
$c = $ws.Range("D23")
$c.Value = "'70.18"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.24%  "
$c = $ws.Range("D24")
$c.Value = "'268.44"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("E25").Value = "  -2.04%  "
$c = $ws.Range("D26")
$c.Value = "'8.24"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +5.13%  "
$c = $ws.Range("D27")
$c.Value = "'26.88"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +2.81%  "
$ws.Range("E28").Value = "  +2.45%  "
$c = $ws.Range("D29")
$c.Value = "'7.24"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -3.35%  "
$ws.Range("E30").Value = "  +0.03%  "
$c = $ws.Range("D31")
$c.Value = "'0.108"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -2.24%  "
$ws.Range("E32").Value = "  -0.91%  "
$c = $ws.Range("D33")
$c.Value = "'34.65"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.59%  "
$ws.Range("E34").Value = "  -0.12%  "
$c = $ws.Range("D35")
$c.Value = "'50.40"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -2.16%  "
$c = $ws.Range("D36")
$c.Value = "'0.0448"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.96%  "
$ws.Range("E37").Value = "  -0.11%  "
$c = $ws.Range("D38")
$c.Value = "'3.34"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +2.38%  "
$ws.Range("E39").Value = "  +7.30%  "
$c = $ws.Range("D40")
$c.Value = "'16.99"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +1.03%  "
$ws.Range("E41").Value = "  +1.06%  "
$ws.Range("E42").Value = "  -0.59%  "
$ws.Range("E43").Value = "  -1.32%  "
$c = $ws.Range("D44")
$c.Value = "'125.12"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("E45").Value = "  +2.41%  "
$c = $ws.Range("D46")
$c.Value = "'21.88"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.56%  "
$ws.Range("E47").Value = "  +3.47%  "
$c = $ws.Range("D48")
$c.Value = "'2.42"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +2.22%  "
$c = $ws.Range("D49")
$c.Value = "'2.030.17"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.54%  "
$c = $ws.Range("D50")
$c.Value = "'3.353.00"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +2.25%  "
$c = $ws.Range("D51")
$c.Value = "'0.0320"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -4.15%  "
